$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = 45436
$ws.Range("D29").Value = 112.4
$ws.Range("D30").Value = 187.2
